$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "69.356.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "3.673.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "685.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "159.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -2.11%  "
$ws.Range("E7").Value2 = "  +0.01%  "
$ws.Range("E8").Value2 = "  -1.10%  "
$ws.Range("E9").Value2 = "  -1.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "7.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -2.17%  "
$ws.Range("E11").Value2 = "  -3.68%  "
$ws.Range("E12").Value2 = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "4.293.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  -0.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "32.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -4.02%  "
$ws.Range("B15").Value2 = "WrappedBTC"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "69.342.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +0.12%  "
$ws.Range("B16").Value2 = "WrappedEther"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "3.663.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -0.50%  "
$ws.Range("E17").Value2 = "  +1.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "15.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -3.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "6.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -4.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "469.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -2.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "9.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.649"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -2.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "79.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "3.820.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -0.15%  "
$ws.Range("E25").Value2 = "  +0.07%  "
$ws.Range("E26").Value2 = "  -3.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "10.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -5.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "9.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -3.71%  "
$ws.Range("E29").Value2 = "  -1.40%  "
$ws.Range("E30").Value2 = "  -5.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "6.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -2.93%  "
$ws.Range("E32").Value2 = "  -0.01%  "
$ws.Range("E33").Value2 = "  -5.93%  "
$ws.Range("E34").Value2 = "  -0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "3.647.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.160"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -1.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "8.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  -4.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "6.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +1.13%  "
$ws.Range("E39").Value2 = "  +0.01%  "
$ws.Range("E40").Value2 = "  +1.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.0897"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -5.17%  "
$ws.Range("E42").Value2 = "  -0.03%  "
$ws.Range("E43").Value2 = "  -1.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "166.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +5.59%  "
$ws.Range("E45").Value2 = "  -1.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.000281"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +0.98%  "
$ws.Range("E47").Value2 = "  -2.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +4.98%  "
$ws.Range("E49").Value2 = "  -1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "27.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "7.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -3.90%  "
